$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Export Worksheet")
$ws1.Columns("A").HorizontalAlignment = -4108
$ws1.Columns("B").HorizontalAlignment = -4108
$ws1.Columns("D").HorizontalAlignment = -4108
